$wb = $excel.ActiveWorkbook

# --- Sheet "soil": rename "con" chamber labels to "ctr" ---
$soil = $wb.Worksheets.Item(2)
$soil.Range("A3").Value = "chs_ctr4"
$soil.Range("A5").Value = "chs_ctr6"

# --- Sheet "branch": replace numeric chamber codes with named labels ---
$branch = $wb.Worksheets.Item(1)
$branch.Range("A2").Value = "chb_irr1"
$branch.Range("A3").Value = "chb_ctr2"
$branch.Range("A4").Value = "chb_irr1"
$branch.Range("A5").Value = "chb_ctr2"
$branch.Range("A6").Value = "chb_irr1"
$branch.Range("A7").Value = "chb_ctr2"

# --- Update selections / active sheet to match final UI state ---
$soil.Select() | Out-Null
$soil.Range("A5").Select() | Out-Null

$branch.Select() | Out-Null
$branch.Range("A6:A7").Select() | Out-Null
